$d = $word.ActiveDocument

# --- Split 1: paragraph "{m:v}" -> split run "{m" into "{" and "m" -----
# Locate the paragraph that starts with "{m:v}" and is NOT inside the table
# (the table cell copy of "{m:v}" must stay untouched).
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("{m:v}") -and ($p.Range.Information(12) -eq $false)) {
        $target1 = $p.Range
    }
}
if ($target1 -eq $null) {
    throw "Could not find target paragraph '{m:v}' outside of a table"
}

$splitPoint1 = $target1.Start + 1  # right after the "{" character, inside the "{m" run
$splitRange1 = $d.Range($splitPoint1, $splitPoint1)
$d.Bookmarks.Add("m2docSplitMark1", $splitRange1) | Out-Null
$d.Bookmarks("m2docSplitMark1").Delete()

# --- Split 2: paragraph "{m:endfor}" -> split run "{m:" into "{" and "m:" --
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("{m:endfor}")) {
        $target2 = $p.Range
    }
}
if ($target2 -eq $null) {
    throw "Could not find target paragraph '{m:endfor}'"
}

$splitPoint2 = $target2.Start + 1  # right after the "{" character, inside the "{m:" run
$splitRange2 = $d.Range($splitPoint2, $splitPoint2)
$d.Bookmarks.Add("m2docSplitMark2", $splitRange2) | Out-Null
$d.Bookmarks("m2docSplitMark2").Delete()
